# Slide 9 ("Fazit"), content placeholder: the bullet that reads
#   "Useless use of cat -> cat crime.csv | wc –l -> besser: wc -l < crime.csv"
# keeps the exact same wording, but the author re-split a couple of runs
# (so "cat" lines up with the rest of the shell snippet instead of "of")
# and switched the shell-command fragments over to the monospace
# "Consolas" font - fixing a sloppy earlier edit ("fixed write error ;)").

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(9)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Anchor on the start of the sentence so we don't depend on hard-coded
# absolute character offsets into the whole placeholder.
$anchor = $tr.Find("Useless use of cat", 0)
$base = $anchor.Start

# --- Re-split "use" / " of " / "cat" -------------------------------------
# Before: "use"[err] " "[plain] "of"[err] " "[plain] "cat"[err] " -> "[plain]
# After : "use"[err] " of "[plain]        "cat"[err] " "[plain] "->"[Consolas] " "[plain]
# Re-stamping each span with its own (unchanged) text forces PowerPoint to
# re-group the runs at these exact boundaries while leaving the rendered
# text untouched.
$tr.Characters($base + 11, 4).Text = " of "
$tr.Characters($base + 15, 3).Text = "cat"

# --- Apply the Consolas font to the shell-command fragments ---------------
$tr.Characters($base + 19, 2).Font.Name  = "Consolas"   # "->"
$tr.Characters($base + 22, 16).Font.Name = "Consolas"   # "cat crime.csv | "
$tr.Characters($base + 38, 2).Font.Name  = "Consolas"   # "wc"
$tr.Characters($base + 40, 7).Font.Name  = "Consolas"   # " –l -> "
$tr.Characters($base + 55, 2).Font.Name  = "Consolas"   # "wc"
$tr.Characters($base + 57, 15).Font.Name = "Consolas"   # " -l < crime.csv"
